$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("knnts")

# Update header row (row 1). Column B ("Arreglo aleatorio óptimo knnts") is
# removed, and everything after it shifts one column to the left.
$ws.Range("A1").Value = "K óptimo"
$ws.Range("B1").Value = "MAE knnts"
$ws.Range("C1").Value = "MSE knnts"
$ws.Range("D1").Value = "RMSE knnts"
$ws.Range("E1").Value = "R2 knnts"

# Update data row (row 2) with the new results.
$ws.Range("A2").Value = 23
$ws.Range("B2").Value = 0.5753903495649837
$ws.Range("C2").Value = 0.6224774416917598
$ws.Range("D2").Value = 0.7889723960264768
$ws.Range("E2").Value = -0.03071067376510817

# The old column F is no longer part of the table; clear it so the used
# range (and dimension) shrinks back down to A1:E2.
$ws.Range("F1:F2").Clear()
